# Fix OCR-extracted values in the output table (corrections from PaddleOCR re-run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "1.000 GHz"     # was "1,000 GHz"
$ws.Range("C1").Value = "1075.615 K"    # was "1075.615 kK]"
$ws.Range("B2").Value = "6.291 dB"      # was "6.23"
$ws.Range("B3").Value = "6.787 dB"      # was "6.787"
$ws.Range("B4").Value = "6.737 d8"      # was "6.737 dB"
$ws.Range("C4").Value = "1077.982 K"    # was "K"
$ws.Range("C5").Value = "1184.725 K"    # was "1184.726 kK"

# Column C's width auto-tracks its (now shorter) best-fit content.
$ws.Columns.Item(3).ColumnWidth = 9.4
